# Update automàtic: dades i banners [2026-02-17 21:50]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-17 21:48:58"
$ws.Range("H2").Value = "'58%"
$ws.Range("E3").Value = "2026-02-17 21:49:01"
$ws.Range("E4").Value = "2026-02-17 21:49:04"
$ws.Range("E5").Value = "2026-02-17 21:49:06"
$ws.Range("E6").Value = "2026-02-17 21:49:09"
$ws.Range("E7").Value = "2026-02-17 21:49:12"
$ws.Range("H7").Value = "'63%"
$ws.Range("E8").Value = "2026-02-17 21:49:14"
$ws.Range("E9").Value = "2026-02-17 21:49:17"
$ws.Range("H9").Value = "'60%"
$ws.Range("N9").Value = "7.5 °C 21:29 TU"
$ws.Range("O9").Value = "12.3 °C"
$ws.Range("E10").Value = "2026-02-17 21:49:20"
$ws.Range("O10").Value = "10.4 °C"
$ws.Range("E11").Value = "2026-02-17 21:49:22"
$ws.Range("H11").Value = "'53%"
$ws.Range("E12").Value = "2026-02-17 21:49:25"
$ws.Range("H12").Value = "'62%"
$ws.Range("N12").Value = "8.7 °C 21:18 TU"
$ws.Range("O12").Value = "12.5 °C"
$ws.Range("E13").Value = "2026-02-17 21:49:27"
$ws.Range("J13").Value = "1018.1 hPa"
$ws.Range("O13").Value = "6.7 °C"
$ws.Range("E14").Value = "2026-02-17 21:49:30"
$ws.Range("H14").Value = "'70%"
$ws.Range("O14").Value = "13.4 °C"
$ws.Range("E15").Value = "2026-02-17 21:49:32"
$ws.Range("H15").Value = "'59%"
$ws.Range("N15").Value = "6.7 °C 21:05 TU"
$ws.Range("O15").Value = "11.9 °C"
$ws.Range("E16").Value = "2026-02-17 21:49:35"
$ws.Range("H16").Value = "'65%"
$ws.Range("M16").Value = "0.0 °C 21:22 TU"
$ws.Range("O16").Value = "-3.3 °C"
$ws.Range("E17").Value = "2026-02-17 21:49:38"
$ws.Range("E18").Value = "2026-02-17 21:49:40"
$ws.Range("J18").Value = "1018.8 hPa"
$ws.Range("E19").Value = "2026-02-17 21:49:43"
$ws.Range("H19").Value = "'75%"
$ws.Range("E20").Value = "2026-02-17 21:49:46"
$ws.Range("H20").Value = "'66%"
$ws.Range("E21").Value = "2026-02-17 21:49:48"
$ws.Range("H21").Value = "'40%"
$ws.Range("O21").Value = "9.6 °C"
$ws.Range("E22").Value = "2026-02-17 21:49:51"
$ws.Range("E23").Value = "2026-02-17 21:49:53"
$ws.Range("K23").Value = "12.2 MJ/m2"
$ws.Range("E24").Value = "2026-02-17 21:49:56"
$ws.Range("E25").Value = "2026-02-17 21:49:59"
$ws.Range("E26").Value = "2026-02-17 21:50:01"
$ws.Range("E27").Value = "2026-02-17 21:50:04"
$ws.Range("H27").Value = "'55%"
$ws.Range("E28").Value = "2026-02-17 21:50:06"
$ws.Range("H28").Value = "'80%"
$ws.Range("J28").Value = "1018.4 hPa"
$ws.Range("E29").Value = "2026-02-17 21:50:09"
$ws.Range("H29").Value = "'67%"
$ws.Range("O29").Value = "11.9 °C"
$ws.Range("E30").Value = "2026-02-17 21:50:12"
$ws.Range("J30").Value = "1018.5 hPa"
$ws.Range("N30").Value = "8.0 °C 21:21 TU"
$ws.Range("O30").Value = "11.2 °C"
$ws.Range("E31").Value = "2026-02-17 21:50:14"
$ws.Range("H31").Value = "'68%"
$ws.Range("E32").Value = "2026-02-17 21:50:16"
$ws.Range("K32").Value = "11.0 MJ/m2"
$ws.Range("E33").Value = "2026-02-17 21:50:19"
$ws.Range("H33").Value = "'44%"
$ws.Range("J33").Value = "1017.5 hPa"
$ws.Range("N33").Value = "3.6 °C 21:28 TU"
$ws.Range("O33").Value = "6.4 °C"
$ws.Range("E34").Value = "2026-02-17 21:50:22"
$ws.Range("H34").Value = "'55%"
$ws.Range("E35").Value = "2026-02-17 21:50:24"
$ws.Range("E36").Value = "2026-02-17 21:50:27"
$ws.Range("J36").Value = "1018.7 hPa"
$ws.Range("N36").Value = "9.7 °C 21:08 TU"
$ws.Range("E37").Value = "2026-02-17 21:50:29"
$ws.Range("J37").Value = "1019.3 hPa"
$ws.Range("N37").Value = "3.9 °C 21:29 TU"
$ws.Range("O37").Value = "7.3 °C"
$ws.Range("E38").Value = "2026-02-17 21:50:32"
$ws.Range("H38").Value = "'78%"
$ws.Range("O38").Value = "10.9 °C"
$ws.Range("E39").Value = "2026-02-17 21:50:34"
$ws.Range("H39").Value = "'60%"
$ws.Range("E40").Value = "2026-02-17 21:50:37"
$ws.Range("H40").Value = "'53%"
$ws.Range("E41").Value = "2026-02-17 21:50:39"
$ws.Range("J41").Value = "1018.2 hPa"
$ws.Range("K41").Value = "12.2 MJ/m2"
$ws.Range("O41").Value = "16.3 °C"
$ws.Range("E42").Value = "2026-02-17 21:50:41"
$ws.Range("H42").Value = "'61%"
$ws.Range("N42").Value = "9.3 °C 21:00 TU"
$ws.Range("O42").Value = "12.7 °C"
$ws.Range("E43").Value = "2026-02-17 21:50:44"
$ws.Range("H43").Value = "'84%"
$ws.Range("E44").Value = "2026-02-17 21:50:47"
$ws.Range("E45").Value = "2026-02-17 21:50:50"
$ws.Range("H45").Value = "'71%"
$ws.Range("N45").Value = "1.8 °C 21:28 TU"
$ws.Range("O45").Value = "5.3 °C"
$ws.Range("E46").Value = "2026-02-17 21:50:52"
$ws.Range("H46").Value = "'58%"
$ws.Range("J46").Value = "1019.5 hPa"
$ws.Range("K46").Value = "14.2 MJ/m2"
$ws.Range("N46").Value = "10.3 °C 21:00 TU"
$ws.Range("O46").Value = "15.2 °C"
